# Regenerated staging template after Meerkat DB changes: a new
# "BusinessKey" column is inserted as the first data column on the
# header row, pushing the existing columns (EndDateID, ID,
# ReportingPeriod, StartDateID, Summary, YearNumber) one slot to the
# right, e.g. A2:F2 -> B2:G2, with A2 becoming "BusinessKey".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing header row (row 2) one column to the right, working
# from the rightmost column back towards A so no value is clobbered
# before it has been copied onward. (Reads use the explicit Value()
# call form - plain ".Value" property access on this host returns the
# COM property descriptor rather than invoking the getter.)
$ws.Range("G2").Value = $ws.Range("F2").Value()
$ws.Range("F2").Value = $ws.Range("E2").Value()
$ws.Range("E2").Value = $ws.Range("D2").Value()
$ws.Range("D2").Value = $ws.Range("C2").Value()
$ws.Range("C2").Value = $ws.Range("B2").Value()
$ws.Range("B2").Value = $ws.Range("A2").Value()

# New first column header.
$ws.Range("A2").Value = "BusinessKey"
